# "prepare the seeds and parse data form csv file function"
# - Users sheet: rename the "name" header (A1) to "username"
# - Users sheet: add a helper column D (value 1) for each existing user row,
#   used as a flag/marker when parsing the CSV seed data
# - Make "Users" the active sheet/tab instead of "Purchases"

$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")

# Rename header of column A on the Users sheet from "name" to "username"
$wsUsers.Range("A1").Value = "username"

# Add a new column D populated with 1 for each existing user row (rows 2-12)
for ($r = 2; $r -le 12; $r++) {
    $wsUsers.Cells.Item($r, 4).Value = 1
}

# Switch the active/selected tab from Purchases to Users
$wsUsers.Activate()

$wb.Save()
